$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected in the source file; unprotect to allow edits, then
# re-protect once all writes are done.
$ws.Unprotect()

# Update the confidential footer note's "as of" date (2021-05-10 -> 2021-05-11).
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns for rows 2-39.
$ws.Cells.Item(2, 4).Value = 0.05767452158711807
$ws.Cells.Item(2, 5).Value = -0.00741032715806067
$ws.Cells.Item(3, 4).Value = 0.05227191334171964
$ws.Cells.Item(3, 5).Value = -0.00384335302208938
$ws.Cells.Item(4, 4).Value = 0.3110134434487704
$ws.Cells.Item(4, 5).Value = -0.0242741551642075
$ws.Cells.Item(5, 4).Value = 0.03373513569010905
$ws.Cells.Item(5, 5).Value = 0.01047488003410146
$ws.Cells.Item(6, 4).Value = 0.03106127019818938
$ws.Cells.Item(6, 5).Value = -0.002168429437536035
$ws.Cells.Item(7, 4).Value = 0.03068431945164187
$ws.Cells.Item(7, 5).Value = -0.01662324773601298
$ws.Cells.Item(8, 4).Value = 0.02880601564752684
$ws.Cells.Item(8, 5).Value = -0.00816350502143659
$ws.Cells.Item(9, 4).Value = 0.02382371012794226
$ws.Cells.Item(9, 5).Value = -0.009018605311745365
$ws.Cells.Item(10, 4).Value = 0.02423217036185897
$ws.Cells.Item(10, 5).Value = -0.009464383113341301
$ws.Cells.Item(11, 4).Value = 0.02264654528605909
$ws.Cells.Item(11, 5).Value = 0.001830244795241409
$ws.Cells.Item(12, 4).Value = 0.02312035069847968
$ws.Cells.Item(12, 5).Value = -0.01141498216409043
$ws.Cells.Item(13, 4).Value = 0.0214359676991804
$ws.Cells.Item(13, 5).Value = -0.0142702116115031
$ws.Cells.Item(14, 4).Value = 0.02183512311803506
$ws.Cells.Item(14, 5).Value = -0.004939347715551379
$ws.Cells.Item(15, 4).Value = 0.02164130804975411
$ws.Cells.Item(15, 5).Value = -0.03066369606003749
$ws.Cells.Item(16, 4).Value = 0.02225172506513106
$ws.Cells.Item(16, 5).Value = -0.02183468364655838
$ws.Cells.Item(17, 4).Value = 0.01936776107857186
$ws.Cells.Item(17, 5).Value = 0.003466724900365925
$ws.Cells.Item(18, 4).Value = 0.01420337724606869
$ws.Cells.Item(18, 5).Value = -0.01661604430945141
$ws.Cells.Item(19, 4).Value = 0.01702654272512022
$ws.Cells.Item(19, 5).Value = -0.0001738828029907991
$ws.Cells.Item(20, 4).Value = 0.01573359637763549
$ws.Cells.Item(20, 5).Value = -0.01646505376344087
$ws.Cells.Item(21, 4).Value = 0.01720418010357739
$ws.Cells.Item(21, 5).Value = -0.03179929689996785
$ws.Cells.Item(22, 4).Value = 0.0133025019696073
$ws.Cells.Item(22, 5).Value = -0.01882233244308773
$ws.Cells.Item(23, 4).Value = 0.01509558212667681
$ws.Cells.Item(23, 5).Value = -0.0107448552176288
$ws.Cells.Item(24, 4).Value = 0.01483578746068319
$ws.Cells.Item(24, 5).Value = -0.01133925835121075
$ws.Cells.Item(25, 4).Value = 0.01405238547438009
$ws.Cells.Item(25, 5).Value = -0.006207674943566555
$ws.Cells.Item(26, 4).Value = 0.01383393379349168
$ws.Cells.Item(26, 5).Value = -0.004624180258954058
$ws.Cells.Item(27, 4).Value = 0.01298021701154928
$ws.Cells.Item(27, 5).Value = 0.005278592375366431
$ws.Cells.Item(28, 4).Value = 0.01390266254110904
$ws.Cells.Item(28, 5).Value = -0.02619330108606355
$ws.Cells.Item(29, 4).Value = 0.01432983814168463
$ws.Cells.Item(29, 5).Value = -0.01279478173607618
$ws.Cells.Item(30, 4).Value = 0.01355965322223102
$ws.Cells.Item(30, 5).Value = -0.03069245165315049
$ws.Cells.Item(31, 4).Value = 0.01241093150820795
$ws.Cells.Item(31, 5).Value = -0.007292802617230087
$ws.Cells.Item(32, 4).Value = 0.01351757008138224
$ws.Cells.Item(32, 5).Value = -0.01144381345723633
$ws.Cells.Item(33, 4).Value = 0.01239813738749765
$ws.Cells.Item(33, 5).Value = -0.009082768325444635
$ws.Cells.Item(34, 4).Value = 0.006033643885060579
$ws.Cells.Item(34, 5).Value = 0.00283896745702128
$ws.Cells.Item(35, 4).Value = 0.00514609141198348
$ws.Cells.Item(35, 5).Value = 0.01723889950481827
$ws.Cells.Item(36, 4).Value = 0.005152118394466849
$ws.Cells.Item(36, 5).Value = 0.02039978656158925
$ws.Cells.Item(37, 4).Value = 0.00506879800504765
$ws.Cells.Item(37, 5).Value = 0.01211982143602164
$ws.Cells.Item(38, 4).Value = 0.004611170282450833
$ws.Cells.Item(38, 5).Value = 0.009103416647557871
$ws.Cells.Item(39, 4).Value = 0.9999999999999998
$ws.Cells.Item(39, 5).Value = -0.01351556108722118

$ws.Protect()
